$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new sheet "2022-Q1" right before the "总计" sheet.
#    Copy an existing quarter sheet as a template so the new sheet keeps
#    the same sheetPr / pageMargins / header style as its siblings, then
#    overwrite its contents with the 2022-Q1 numbers.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($wb.Worksheets.Item("总计"))

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Numeric-looking text (fund codes / decimal figures) must be forced to
# the Text format first, otherwise Excel silently coerces them to real
# numbers (and "000906" would lose its leading zeros). ClearFormats()
# afterwards drops the leftover "@" number-format style so the cell
# keeps the plain default style, matching the source rows.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "270023"
$newSheet.Range("C2").Value = "广发全球精选股票(QDII)"
$newSheet.Range("D2").Value = "25.53"
$newSheet.Range("E2").Value = "78.43"
$newSheet.Range("F2").Value = "4.06"
$newSheet.Range("G2").Value = "1.0365"
$newSheet.Range("B2:G2").ClearFormats()
$newSheet.Range("H2").Value = 6

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3:G3").NumberFormat = "@"
$newSheet.Range("B3").Value = "000906"
$newSheet.Range("C3").Value = "广发全球精选股票(QDII)美元现汇"
$newSheet.Range("D3").Value = "25.53"
$newSheet.Range("E3").Value = "78.43"
$newSheet.Range("F3").Value = "4.06"
$newSheet.Range("G3").Value = "1.0365"
$newSheet.Range("B3:G3").ClearFormats()
$newSheet.Range("H3").Value = 6

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing rows down by one and renumbering the leading index column.
#    Re-fetch the sheet by name (the reference captured before the sheet
#    copy/insert above does not track the sheet correctly afterwards).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The inserted row picked up formatting from the row below it; clear the
# data cells back to the default (unstyled) look, then copy the bordered
# "index" style from the row beneath onto the new index cell.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 2.07

# Renumber the rest of the index column (0,1,2,3...) now that everything
# shifted down by one row.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
